# New Local_Ks, working on dynamics from spreadsheets
#
# Update the NAVARCH sheet's local-calculation inputs:
#   - B23 (Vehicles weight input) cleared out
#   - A39 (new local Ks source value) updated
#   - C39 / D39 (x, z values tied to the old Ks row) cleared out
# Everything else on the sheet (B24/C24/D24, B25:D25, B27:D27, D29, B31,
# B32, B33, B39) recalculates automatically from these inputs.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NAVARCH")

# Clear the old "Vehicles" weight figure in the summary block.
$ws.Range("B23").ClearContents()

# New Local_Ks source row: update A39, and clear the now-unused C39/D39.
$ws.Range("A39").Value = 778.52111258704656
$ws.Range("C39").ClearContents()
$ws.Range("D39").ClearContents()
